$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting the existing D:K data block to E:L.
$ws.Range("D:D").EntireColumn.Insert()

# Copy formatting (number format / font / alignment) from the shifted column E
# into the newly inserted column D for every row that carries a data block,
# skipping the section-header-only rows which have no D:L cells at all.
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Match the column width of the newly inserted column to its neighbour
$ws.Range("D1").ColumnWidth = $ws.Range("E1").ColumnWidth

# Populate the new column D with the latest (newest) fiscal-period figures
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 2384200
$ws.Range("D9").Value = 1808100
$ws.Range("D10").Value = 576200
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 2175300
$ws.Range("D18").Value = 209000
$ws.Range("D20").Value = 600
$ws.Range("D21").Value = 248900
$ws.Range("D22").Value = 28700
$ws.Range("D23").Value = 180800
$ws.Range("D24").Value = 46100
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 134800
$ws.Range("D27").Value = 134800
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -600
$ws.Range("D33").Value = 134800
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 134800
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 100900
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 407100
$ws.Range("D44").Value = 169000
$ws.Range("D45").Value = 27700
$ws.Range("D46").Value = 704700
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 168000
$ws.Range("D49").Value = 1563400
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 18500
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 2454500
$ws.Range("D57").Value = 313200
$ws.Range("D58").Value = 26900
$ws.Range("D59").Value = 104200
$ws.Range("D60").Value = 444300
$ws.Range("D61").Value = 716600
$ws.Range("D62").Value = 221600
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 1382400
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 441900
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 1072100
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 134800
$ws.Range("D83").Value = 39400
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 167200
$ws.Range("D91").Value = -52500
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -551800
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 429100
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 44400
